# people_widget: change the "people" field from select_or_add_multiple to a
# textarea, and switch the example value from a comma-separated string to a
# newline-separated one (one name per line), matching how a textarea widget
# expects multi-line input.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# C3 = type column for the "people" row
$ws.Range("C3").Value = "textarea"

# M3 = example value column - now newline separated instead of comma separated
$nl = [char]10
$ws.Range("M3").Value = '"Ruben' + $nl + 'Jessica"'
$ws.Range("M3").WrapText = $true

# Move the active selection/cell to M3 to match the updated workbook state
[void]$ws.Activate()
[void]$ws.Range("M3").Select()
